# semana 24 de 2025
# Adds two new weekly columns (Z = week 23, AA = week 24) to the weekly
# IRA-hospital report, mirroring the header style of the existing week
# columns and filling in the per-establishment counts for the two new weeks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new week numbers, formatted like the existing ones ---
# (the existing week headers are stored as text, so a leading apostrophe
# keeps "23"/"24" as text instead of being auto-converted to numbers;
# re-pasting Y1's format afterwards matches the bold/centered header style
# without leaving a stray "quote prefix" style variant behind)
$ws.Range("Z1").Value = "'23"
$ws.Range("AA1").Value = "'24"
$ws.Range("Y1").Copy()
$ws.Range("Z1:AA1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows: counts for week 23 (Z) and week 24 (AA) ---
$weekData = @{
    2  = @{ Z = 0;  AA = 0 }
    3  = @{ Z = 0 }
    4  = @{ Z = 0;  AA = 0 }
    5  = @{ Z = 0;  AA = 0 }
    6  = @{ Z = 32; AA = 21 }
    7  = @{ Z = 2;  AA = 6 }
    8  = @{ Z = 2;  AA = 4 }
    9  = @{ Z = 0;  AA = 0 }
    10 = @{ Z = 0 }
    11 = @{ Z = 0 }
    12 = @{ Z = 0;  AA = 0 }
    13 = @{ Z = 0;  AA = 0 }
    14 = @{ Z = 0;  AA = 0 }
    15 = @{ Z = 0;  AA = 0 }
    16 = @{ Z = 0;  AA = 0 }
    17 = @{ Z = 0;  AA = 0 }
    18 = @{ Z = 0 }
    19 = @{ Z = 0 }
    21 = @{ Z = 0;  AA = 0 }
    22 = @{ Z = 0;  AA = 0 }
    23 = @{ AA = 0 }
    24 = @{ Z = 1;  AA = 4 }
    25 = @{ Z = 0;  AA = 0 }
    27 = @{ Z = 0;  AA = 0 }
    28 = @{ Y = 2;  Z = 3;  AA = 3 }
    29 = @{ Y = 17; Z = 16; AA = 10 }
    30 = @{ Z = 0;  AA = 0 }
    31 = @{ Z = 0;  AA = 0 }
    33 = @{ Z = 0;  AA = 0 }
    34 = @{ Z = 28; AA = 26 }
    35 = @{ Y = 2;  Z = 1;  AA = 0 }
    36 = @{ Z = 0;  AA = 0 }
    37 = @{ Z = 0;  AA = 0 }
    39 = @{ Z = 0;  AA = 0 }
    40 = @{ Z = 0 }
    41 = @{ Z = 0 }
    42 = @{ Z = 0;  AA = 0 }
    43 = @{ Z = 0;  AA = 0 }
    44 = @{ Z = 0;  AA = 0 }
    45 = @{ Z = 0;  AA = 0 }
    46 = @{ Z = 0;  AA = 0 }
    47 = @{ Z = 0;  AA = 0 }
    48 = @{ Z = 0;  AA = 0 }
    49 = @{ AA = 0 }
    50 = @{ Y = 0;  Z = 0 }
    51 = @{ Z = 0;  AA = 0 }
    52 = @{ Z = 0;  AA = 0 }
    53 = @{ Z = 0;  AA = 0 }
    54 = @{ Z = 0;  AA = 0 }
    55 = @{ Z = 0;  AA = 0 }
    56 = @{ Z = 0;  AA = 0 }
}

foreach ($row in $weekData.Keys) {
    $cells = $weekData[$row]
    foreach ($col in $cells.Keys) {
        $ws.Range("$col$row").Value = $cells[$col]
    }
}
